$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Cells.Item(2, 4) '27.081.91'
Set-TextValue $ws.Cells.Item(2, 5) '  +0.50%  '

# Row 3
Set-TextValue $ws.Cells.Item(3, 4) '1.680.19'
Set-TextValue $ws.Cells.Item(3, 5) '  +0.50%  '

# Row 4
Set-TextValue $ws.Cells.Item(4, 5) '  +0.02%  '

# Row 5
Set-TextValue $ws.Cells.Item(5, 4) '215.02'
Set-TextValue $ws.Cells.Item(5, 5) '  +0.13%  '

# Row 6
Set-TextValue $ws.Cells.Item(6, 5) '  +0.06%  '

# Row 8
Set-TextValue $ws.Cells.Item(8, 5) '  +1.97%  '

# Row 9
Set-TextValue $ws.Cells.Item(9, 4) '21.29'
Set-TextValue $ws.Cells.Item(9, 5) '  +5.77%  '

# Row 10
Set-TextValue $ws.Cells.Item(10, 5) '  +0.49%  '

# Row 11
Set-TextValue $ws.Cells.Item(11, 4) '0.0884'
Set-TextValue $ws.Cells.Item(11, 5) '  -0.62%  '

# Row 12
Set-TextValue $ws.Cells.Item(12, 4) '1.917.45'
Set-TextValue $ws.Cells.Item(12, 5) '  +0.52%  '

# Row 13
Set-TextValue $ws.Cells.Item(13, 4) '1.675.03'
Set-TextValue $ws.Cells.Item(13, 5) '  +0.50%  '

# Row 14
Set-TextValue $ws.Cells.Item(14, 5) '  +1.25%  '

# Row 15
Set-TextValue $ws.Cells.Item(15, 5) '  +2.10%  '

# Row 16
Set-TextValue $ws.Cells.Item(16, 4) '66.12'
Set-TextValue $ws.Cells.Item(16, 5) '  +0.79%  '

# Row 17
Set-TextValue $ws.Cells.Item(17, 4) '27.085.74'
Set-TextValue $ws.Cells.Item(17, 5) '  +0.44%  '

# Row 18
Set-TextValue $ws.Cells.Item(18, 4) '237.81'
Set-TextValue $ws.Cells.Item(18, 5) '  +1.45%  '

# Row 19
Set-TextValue $ws.Cells.Item(19, 4) '8.14'
Set-TextValue $ws.Cells.Item(19, 5) '  +0.95%  '

# Row 20
Set-TextValue $ws.Cells.Item(20, 5) '  +1.54%  '

# Row 21
Set-TextValue $ws.Cells.Item(21, 5) '  +0.03%  '

# Row 22
Set-TextValue $ws.Cells.Item(22, 5) '  +1.46%  '

# Row 23
Set-TextValue $ws.Cells.Item(23, 4) '9.38'
Set-TextValue $ws.Cells.Item(23, 5) '  +2.52%  '

# Row 24
Set-TextValue $ws.Cells.Item(24, 5) '  -2.12%  '

# Row 25
Set-TextValue $ws.Cells.Item(25, 4) '146.77'
Set-TextValue $ws.Cells.Item(25, 5) '  +0.74%  '

# Row 26
Set-TextValue $ws.Cells.Item(26, 4) '7.21'
Set-TextValue $ws.Cells.Item(26, 5) '  +1.01%  '

# Row 27
Set-TextValue $ws.Cells.Item(27, 4) '16.32'
Set-TextValue $ws.Cells.Item(27, 5) '  +2.26%  '

# Row 28
Set-TextValue $ws.Cells.Item(28, 5) '  +0.68%  '

# Row 29
Set-TextValue $ws.Cells.Item(29, 4) '1.00'
Set-TextValue $ws.Cells.Item(29, 5) '  +0.13%  '

# Row 30
Set-TextValue $ws.Cells.Item(30, 5) '  +0.23%  '

# Row 31
Set-TextValue $ws.Cells.Item(31, 5) '  +0.14%  '

# Row 32
Set-TextValue $ws.Cells.Item(32, 2) 'Maker'
Set-TextValue $ws.Cells.Item(32, 3) 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue $ws.Cells.Item(32, 4) '1.554.85'
Set-TextValue $ws.Cells.Item(32, 5) '  +5.67%  '

# Row 33
Set-TextValue $ws.Cells.Item(33, 2) 'Filecoin'
Set-TextValue $ws.Cells.Item(33, 3) 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue $ws.Cells.Item(33, 4) '3.36'
Set-TextValue $ws.Cells.Item(33, 5) '  +1.00%  '

# Row 34
Set-TextValue $ws.Cells.Item(34, 5) '  +1.80%  '

# Row 35
Set-TextValue $ws.Cells.Item(35, 5) '  +2.36%  '

# Row 36
Set-TextValue $ws.Cells.Item(36, 4) '0.604'
Set-TextValue $ws.Cells.Item(36, 5) '  +4.77%  '

# Row 37
Set-TextValue $ws.Cells.Item(37, 5) '  +4.91%  '

# Row 38
Set-TextValue $ws.Cells.Item(38, 5) '  -1.21%  '

# Row 39
Set-TextValue $ws.Cells.Item(39, 5) '  +1.97%  '

# Row 40
Set-TextValue $ws.Cells.Item(40, 5) '  +0.97%  '

# Row 41
Set-TextValue $ws.Cells.Item(41, 5) '  +0.00%  '

# Row 42
Set-TextValue $ws.Cells.Item(42, 4) '68.61'
Set-TextValue $ws.Cells.Item(42, 5) '  +2.90%  '

# Row 43
Set-TextValue $ws.Cells.Item(43, 4) '5.63'
Set-TextValue $ws.Cells.Item(43, 5) '  -2.69%  '

# Row 44
Set-TextValue $ws.Cells.Item(44, 5) '  -1.75%  '

# Row 45
Set-TextValue $ws.Cells.Item(45, 4) '1.824.79'
Set-TextValue $ws.Cells.Item(45, 5) '  +0.64%  '

# Row 46
Set-TextValue $ws.Cells.Item(46, 5) '  +0.52%  '

# Row 47
Set-TextValue $ws.Cells.Item(47, 4) '90.76'
Set-TextValue $ws.Cells.Item(47, 5) '  +0.35%  '

# Row 48
Set-TextValue $ws.Cells.Item(48, 2) 'RenderToken'
Set-TextValue $ws.Cells.Item(48, 3) 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws.Cells.Item(48, 4) '1.57'
Set-TextValue $ws.Cells.Item(48, 5) '  +2.65%  '

# Row 49
Set-TextValue $ws.Cells.Item(49, 2) 'BabyDogeCoin'
Set-TextValue $ws.Cells.Item(49, 3) 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-TextValue $ws.Cells.Item(49, 4) '0.0₆0107'
Set-TextValue $ws.Cells.Item(49, 5) '  +1.15%  '

# Row 50
Set-TextValue $ws.Cells.Item(50, 4) '0.104'
Set-TextValue $ws.Cells.Item(50, 5) '  +3.44%  '

# Row 51
Set-TextValue $ws.Cells.Item(51, 4) '8.08'
Set-TextValue $ws.Cells.Item(51, 5) '  +4.95%  '
